{"js": "// Remove the delivery-entry paragraphs that follow the \"Deliveres for ...\"\n// heading, leaving only the heading paragraph in the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Keep the first paragraph (the \"Deliveres for DD/MM/YYYY\" heading) and\n// delete every paragraph that comes after it (the Lewis Luck / Brian Binks /\n// Joseph Jones delivery blocks).\nfor (let i = paragraphs.items.length - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the delivery-entry paragraphs that follow the \"Deliveres for ...\"\n# heading, leaving only the heading paragraph in the document.\n$d = $word.ActiveDocument\n\n# Keep paragraph 1 (the \"Deliveres for DD/MM/YYYY\" heading) and delete every\n# paragraph after it (the Lewis Luck / Brian Binks / Joseph Jones delivery\n# blocks). Iterate backwards so the indices stay valid as we delete.\nfor ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n"}
